# variation_coefficients.xlsx - apply corrected substrate values
# (reordered rows, corrected TS reference) and highlight the updated
# VC values, per commit message:
# "Version with corrected substrate values, number of steady state
#  simulation changed to 500 and disturbance feeding on scenario 1b
#  corrected so OLR is no more than 1.2"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the substrate table (rows 2-6) with the corrected values
# and the new row order: XP, XL, XA, BMP, TS. The TS row's reference
# now also points to Delory (2025) instead of Hafner (2018). ---

$ws.Range("A2").Value = "XP"
$ws.Range("B2").Value = 5.52
$ws.Range("C2").Value = "Delory (2025) (submitted)"

$ws.Range("A3").Value = "XL"
$ws.Range("B3").Value = 10.039999999999999
$ws.Range("C3").Value = "Delory (2025) (submitted)"

$ws.Range("A4").Value = "XA"
$ws.Range("B4").Value = 7.4
$ws.Range("C4").Value = "Delory (2025) (submitted)"

$ws.Range("A5").Value = "BMP"
$ws.Range("B5").Value = 8.23
$ws.Range("C5").Value = "Hafner (2018) (mean of all 4 substrates, using robust mean and SD)"

$ws.Range("A6").Value = "TS"
$ws.Range("B6").Value = 1.94
$ws.Range("C6").Value = "Delory (2025) (submitted)"

# Highlight the corrected VC values (column B) with the light green
# "Accent 6, Lighter 80%" fill (theme accent6 #70AD47 tinted 80% ->
# #E2F0D9) used to flag the updated figures.
$ws.Range("B2:B6").Interior.Color = 14282978

# Update the selected cell shown in the sheet view.
$ws.Range("D9").Select() | Out-Null

# Keep the page in portrait orientation.
$ws.PageSetup.Orientation = 1

$wb.Save()
